$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update (prices + 1h volume %) as of the GitHub Actions run.
# Column D holds price text; some values are unambiguous numeric strings that
# Excel would otherwise auto-convert to a Number cell on assignment, so those
# cells are switched to Text format first to keep them as literal strings.

# Row 2
$ws.Range("D2").Value = "27.087.24"
$ws.Range("E2").Value = "  +8.42%  "

# Row 3
$ws.Range("D3").Value = "1.756.53"
$ws.Range("E3").Value = "  +4.71%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9953"
$ws.Range("E4").Value = "  -0.46%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.45"
$ws.Range("E5").Value = "  +1.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3754"
$ws.Range("E7").Value = "  +2.81%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3421"
$ws.Range("E8").Value = "  +5.48%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.41"
$ws.Range("E9").Value = "  +3.19%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("E10").Value = "  +4.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07592"
$ws.Range("E11").Value = "  +7.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.415"
$ws.Range("E13").Value = "  +5.29%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.66"
$ws.Range("E14").Value = "  +4.96%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.119"
$ws.Range("E15").Value = "  +7.34%  "

# Row 16
$ws.Range("D16").Value = "1.748.55"
$ws.Range("E16").Value = "  +4.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").Value = "  +4.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06740"
$ws.Range("E18").Value = "  +2.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.35"
$ws.Range("E19").Value = "  +5.70%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9956"
$ws.Range("E20").Value = "  -0.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.87"
$ws.Range("E21").Value = "  +6.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.290"
$ws.Range("E22").Value = "  +6.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.98"
$ws.Range("E23").Value = "  +0.86%  "

# Row 24
$ws.Range("D24").Value = "26.977.99"
$ws.Range("E24").Value = "  +8.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.444"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.498"
$ws.Range("E26").Value = "  +26.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.453"
$ws.Range("E27").Value = "  +2.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.39"
$ws.Range("E28").Value = "  +2.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.77"
$ws.Range("E29").Value = "  +5.47%  "

# Row 30
$ws.Range("D30").Value = "1.949.64"
$ws.Range("E30").Value = "  +5.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.95"
$ws.Range("E31").Value = "  +5.93%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.126"
$ws.Range("E32").Value = "  +0.89%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.102"
$ws.Range("E33").Value = "  +5.21%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08667"
$ws.Range("E34").Value = "  +2.66%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.694"
$ws.Range("E35").Value = "  +3.21%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.03"
$ws.Range("E36").Value = "  +5.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.493"
$ws.Range("E37").Value = "  +6.40%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02367"
$ws.Range("E38").Value = "  +5.72%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06348"
$ws.Range("E39").Value = "  +4.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2196"
$ws.Range("E40").Value = "  +5.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.611"
$ws.Range("E41").Value = "  +4.43%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6381"
$ws.Range("E42").Value = "  +7.00%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.228"
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.38"
$ws.Range("E44").Value = "  +4.35%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9973"
$ws.Range("E45").Value = "  -0.14%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6292"
$ws.Range("E46").Value = "  +9.79%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.915"
$ws.Range("E47").Value = "  +1.69%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.30"
$ws.Range("E48").Value = "  +4.52%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.092"
$ws.Range("E49").Value = "  +6.43%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07251"
$ws.Range("E50").Value = "  +3.59%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.86"
$ws.Range("E51").Value = "  +6.01%  "

